# Auto-generated edit script applying the Alpha_Profits.xlsx diff
# Updates leve-profit calculation cells (H..N columns) across all 8 job sheets
# matching the scheduled-runner price refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 433.66666
$ws.Range("I12").Value = 378.33334
$ws.Range("J12").Value = 599.6667
$ws.Range("K12").Value = 378.33334
$ws.Range("L12").Value = 599.6667
$ws.Range("M12").Value = -208.33334
$ws.Range("N12").Value = -939.6667
$ws.Range("H70").Value = 4282.5713
$ws.Range("I70").Value = 2700
$ws.Range("J70").Value = 5469.5
$ws.Range("K70").Value = 8100
$ws.Range("L70").Value = 16408.5
$ws.Range("M70").Value = -7830
$ws.Range("N70").Value = -16948.5
$ws.Range("H73").Value = 4282.5713
$ws.Range("I73").Value = 2700
$ws.Range("J73").Value = 5469.5
$ws.Range("K73").Value = 8100
$ws.Range("L73").Value = 16408.5
$ws.Range("M73").Value = -7164
$ws.Range("N73").Value = -18280.5
$ws.Range("H96").Value = 827.8333
$ws.Range("I96").Value = 819.25
$ws.Range("J96").Value = 845
$ws.Range("K96").Value = 2457.75
$ws.Range("L96").Value = 2535
$ws.Range("M96").Value = -1084.75
$ws.Range("N96").Value = -5281
$ws.Range("H100").Value = 2775.7
$ws.Range("I100").Value = 3502.5
$ws.Range("J100").Value = 2594
$ws.Range("K100").Value = 3502.5
$ws.Range("L100").Value = 2594
$ws.Range("M100").Value = -2961.5
$ws.Range("N100").Value = -3676
$ws.Range("H132").Value = 37538.586
$ws.Range("I132").Value = 38772.32
$ws.Range("K132").Value = 116316.96
$ws.Range("M132").Value = -113786.96

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4704.65
$ws.Range("I32").Value = 4945.0527
$ws.Range("J32").Value = 137
$ws.Range("K32").Value = 4945.0527
$ws.Range("L32").Value = 137
$ws.Range("M32").Value = -4658.0527
$ws.Range("N32").Value = -711
$ws.Range("H45").Value = 1789.8334
$ws.Range("I45").Value = 1763.6875
$ws.Range("K45").Value = 1763.6875
$ws.Range("M45").Value = -1386.6875
$ws.Range("H132").Value = 13893851
$ws.Range("I132").Value = 3521.6785
$ws.Range("J132").Value = 62510000
$ws.Range("K132").Value = 10565.0355
$ws.Range("L132").Value = 187530000
$ws.Range("M132").Value = -8035.0355
$ws.Range("N132").Value = -187535060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3464.6667
$ws.Range("I105").Value = 4249
$ws.Range("K105").Value = 4249
$ws.Range("M105").Value = -2502
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 15722.571
$ws.Range("I6").Value = 5005
$ws.Range("J6").Value = 20009.6
$ws.Range("K6").Value = 5005
$ws.Range("L6").Value = 20009.6
$ws.Range("M6").Value = -4892
$ws.Range("N6").Value = -20235.6
$ws.Range("H26").Value = 11599.8
$ws.Range("I26").Value = 8333
$ws.Range("J26").Value = 16500
$ws.Range("K26").Value = 8333
$ws.Range("L26").Value = 16500
$ws.Range("M26").Value = -8046
$ws.Range("N26").Value = -17074
$ws.Range("H31").Value = 2134.5356
$ws.Range("I31").Value = 1290.6428
$ws.Range("K31").Value = 1290.6428
$ws.Range("M31").Value = -995.6428000000001
$ws.Range("H34").Value = 2134.5356
$ws.Range("I34").Value = 1290.6428
$ws.Range("K34").Value = 1290.6428
$ws.Range("M34").Value = -1088.6428
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H99").Value = 2222.4
$ws.Range("I99").Value = 2202.6667
$ws.Range("J99").Value = 2252
$ws.Range("K99").Value = 2202.6667
$ws.Range("L99").Value = 2252
$ws.Range("M99").Value = -704.6667000000002
$ws.Range("N99").Value = -5248
$ws.Range("H107").Value = 3307.7222
$ws.Range("I107").Value = 1208.75
$ws.Range("K107").Value = 1208.75
$ws.Range("M107").Value = 711.25
$ws.Range("H126").Value = 2222.4
$ws.Range("I126").Value = 2202.6667
$ws.Range("J126").Value = 2252
$ws.Range("K126").Value = 6608.000100000001
$ws.Range("L126").Value = 6756
$ws.Range("M126").Value = -4138.000100000001
$ws.Range("N126").Value = -11696
$ws.Range("H132").Value = 2249.5
$ws.Range("I132").Value = 2249.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6748.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4218.5
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 2213.0557
$ws.Range("I134").Value = 2177.25
$ws.Range("K134").Value = 6531.75
$ws.Range("M134").Value = -3996.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 95957.09
$ws.Range("I4").Value = 506.73334
$ws.Range("K4").Value = 1520.20002
$ws.Range("M4").Value = -1408.20002
$ws.Range("H23").Value = 2585.4348
$ws.Range("I23").Value = 3201.3333
$ws.Range("J23").Value = 2189.5
$ws.Range("K23").Value = 9603.999899999999
$ws.Range("L23").Value = 6568.5
$ws.Range("M23").Value = -9368.999899999999
$ws.Range("N23").Value = -7038.5
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3566
$ws.Range("H34").Value = 6087.8696
$ws.Range("I34").Value = 933.8
$ws.Range("J34").Value = 7519.5557
$ws.Range("K34").Value = 2801.4
$ws.Range("L34").Value = 22558.6671
$ws.Range("M34").Value = -2717.4
$ws.Range("N34").Value = -22726.6671
$ws.Range("H37").Value = 120606.4
$ws.Range("J37").Value = 120606.4
$ws.Range("L37").Value = 361819.2
$ws.Range("N37").Value = -362043.2
$ws.Range("H92").Value = 194.13637
$ws.Range("J92").Value = 193.7
$ws.Range("L92").Value = 581.0999999999999
$ws.Range("N92").Value = -3077.1
$ws.Range("H122").Value = 43192.418
$ws.Range("J122").Value = 57313.11
$ws.Range("L122").Value = 515817.99
$ws.Range("N122").Value = -520717.99
$ws.Range("H129").Value = 943.8
$ws.Range("J129").Value = 804.75
$ws.Range("L129").Value = 2414.25
$ws.Range("N129").Value = -12414.25
$ws.Range("H131").Value = 405283.25
$ws.Range("J131").Value = 424098.75
$ws.Range("L131").Value = 1272296.25
$ws.Range("N131").Value = -1282376.25
$ws.Range("H137").Value = 686002.4399999999
$ws.Range("J137").Value = 837425.4399999999
$ws.Range("L137").Value = 2512276.32
$ws.Range("N137").Value = -2522476.32
$ws.Range("H140").Value = 2886.6365
$ws.Range("I140").Value = 1676.9412
$ws.Range("K140").Value = 5030.8236
$ws.Range("M140").Value = 149.1764000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 3353000
$ws.Range("I20").Value = 5014000
$ws.Range("J20").Value = 31000
$ws.Range("K20").Value = 5014000
$ws.Range("L20").Value = 31000
$ws.Range("M20").Value = -5013755
$ws.Range("N20").Value = -31490
$ws.Range("H70").Value = 4711.4287
$ws.Range("I70").Value = 4603.5
$ws.Range("K70").Value = 4603.5
$ws.Range("M70").Value = -4333.5
$ws.Range("H73").Value = 4711.4287
$ws.Range("I73").Value = 4603.5
$ws.Range("K73").Value = 4603.5
$ws.Range("M73").Value = -3667.5
$ws.Range("H122").Value = 3663.5
$ws.Range("I122").Value = 3663.5
$ws.Range("K122").Value = 10990.5
$ws.Range("M122").Value = -8540.5
$ws.Range("H133").Value = 90780
$ws.Range("J133").Value = 90780
$ws.Range("L133").Value = 90780
$ws.Range("N133").Value = -100900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 7930.2856
$ws.Range("J43").Value = 7930.2856
$ws.Range("L43").Value = 7930.2856
$ws.Range("N43").Value = -8316.285599999999
$ws.Range("H63").Value = 60084.5
$ws.Range("J63").Value = 60084.5
$ws.Range("L63").Value = 60084.5
$ws.Range("N63").Value = -61582.5
$ws.Range("H66").Value = 60084.5
$ws.Range("J66").Value = 60084.5
$ws.Range("L66").Value = 180253.5
$ws.Range("N66").Value = -187741.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 60084
$ws.Range("J50").Value = 60084
$ws.Range("L50").Value = 60084
$ws.Range("N50").Value = -61346
$ws.Range("H69").Value = 38695.168
$ws.Range("J69").Value = 38695.168
$ws.Range("L69").Value = 38695.168
$ws.Range("N69").Value = -40193.168
$ws.Range("H72").Value = 38695.168
$ws.Range("J72").Value = 38695.168
$ws.Range("L72").Value = 116085.504
$ws.Range("N72").Value = -123573.504
$ws.Range("H126").Value = 2676.0588
$ws.Range("I126").Value = 2143.889
$ws.Range("J126").Value = 3274.75
$ws.Range("K126").Value = 6431.667
$ws.Range("L126").Value = 9824.25
$ws.Range("M126").Value = -3961.667
$ws.Range("N126").Value = -14764.25
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
